$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.995.97"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.121.83"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.72"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.44"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.40"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.14"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "3.636.40"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "66.945.11"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.16"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "3.120.31"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.24"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.27"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.87"
$ws.Range("E22").Value = "  +4.80%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.33"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.90"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.59"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "0.0₃0950"
$ws.Range("E33").Value = "  -6.76%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.977"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "46.92"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.17"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "2.822.99"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "382.81"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("E46").Value = "  -9.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.94"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.95"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("E51").Value = "  -0.84%  "
